$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the header style/format from G1 (the
# existing "sum" header) so it matches the rest of the header row.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the new data value in H2 (plain number, same as the other data cells).
$ws.Range("H2").Value = 0
